$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.349.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.181.70'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.38%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.613'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.84%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '69.91'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.33%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.579'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.92'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0924'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.83%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.48'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.21%  '
$ws.Range("E13").Value = '  -1.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.505.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.00'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.98%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.194.77'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.800'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.178.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000100'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.79'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.92'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.67'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '226.80'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.92'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.73%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.54'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.55%  '
$ws.Range("E29").Value = '  -2.24%  '
$ws.Range("E30").Value = '  +0.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.49'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.57'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0768'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.12'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -9.03%  '
$ws.Range("E36").Value = '  -2.81%  '
$ws.Range("E37").Value = '  -7.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.10'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0284'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.82%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.08'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.62'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -11.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '59.39'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -8.42%  '
$ws.Range("E44").Value = '  -2.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0977'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.26'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.10%  '
$ws.Range("E48").Value = '  -1.61%  '
$ws.Range("E49").Value = '  -2.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.77%  '
$ws.Range("E51").Value = '  -2.45%  '
